# Scheduled runner update: refresh Leve profit calculations (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across sheets
# ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR with newly pulled market data.
$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 8627.833000000001
$ws.Range("I80").Value = 3838.3572
$ws.Range("J80").Value = 15333.1
$ws.Range("K80").Value = 11515.0716
$ws.Range("L80").Value = 45999.3
$ws.Range("M80").Value = -10517.0716
$ws.Range("N80").Value = -47995.3
$ws.Range("H83").Value = 8627.833000000001
$ws.Range("I83").Value = 3838.3572
$ws.Range("J83").Value = 15333.1
$ws.Range("K83").Value = 34545.2148
$ws.Range("L83").Value = 137997.9
$ws.Range("M83").Value = -29553.2148
$ws.Range("N83").Value = -147981.9
$ws.Range("H94").Value = 1000
$ws.Range("I94").Value = 1000
$ws.Range("K94").Value = 1000
$ws.Range("M94").Value = -549
$ws.Range("H132").Value = 4525.0977
$ws.Range("I132").Value = 1897.4642
$ws.Range("K132").Value = 5692.392599999999
$ws.Range("M132").Value = -3162.392599999999
$ws.Range("H135").Value = 1622.9412
$ws.Range("I135").Value = 792.44446
$ws.Range("J135").Value = 2557.25
$ws.Range("K135").Value = 7132.00014
$ws.Range("L135").Value = 23015.25
$ws.Range("M135").Value = -4597.00014
$ws.Range("N135").Value = -28085.25

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1145
$ws.Range("I2").Value = 1117.5
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 1117.5
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -1004.5
$ws.Range("N2").Value = -1426
$ws.Range("H45").Value = 1406.4783
$ws.Range("I45").Value = 1196.6
$ws.Range("K45").Value = 1196.6
$ws.Range("M45").Value = -819.5999999999999
$ws.Range("H61").Value = 3669.15
$ws.Range("I61").Value = 3054.9333
$ws.Range("J61").Value = 4037.68
$ws.Range("K61").Value = 3054.9333
$ws.Range("L61").Value = 4037.68
$ws.Range("M61").Value = -2842.9333
$ws.Range("N61").Value = -4461.68
$ws.Range("H74").Value = 2764.484
$ws.Range("I74").Value = 1560.25
$ws.Range("J74").Value = 4049
$ws.Range("K74").Value = 1560.25
$ws.Range("L74").Value = 4049
$ws.Range("M74").Value = -686.25
$ws.Range("N74").Value = -5797
$ws.Range("H77").Value = 2764.484
$ws.Range("I77").Value = 1560.25
$ws.Range("J77").Value = 4049
$ws.Range("K77").Value = 7801.25
$ws.Range("L77").Value = 20245
$ws.Range("M77").Value = -3433.25
$ws.Range("N77").Value = -28981
$ws.Range("H116").Value = 1145
$ws.Range("I116").Value = 1117.5
$ws.Range("J116").Value = 1200
$ws.Range("K116").Value = 1117.5
$ws.Range("L116").Value = 1200
$ws.Range("M116").Value = 1176.5
$ws.Range("N116").Value = -5788
$ws.Range("H132").Value = 3320.5112
$ws.Range("I132").Value = 3245.2068
$ws.Range("J132").Value = 3457
$ws.Range("K132").Value = 9735.6204
$ws.Range("L132").Value = 10371
$ws.Range("M132").Value = -7205.6204
$ws.Range("N132").Value = -15431
$ws.Range("H134").Value = 31962.455
$ws.Range("J134").Value = 31962.455
$ws.Range("L134").Value = 31962.455
$ws.Range("N134").Value = -42102.455
$ws.Range("H136").Value = 3669.15
$ws.Range("I136").Value = 3054.9333
$ws.Range("J136").Value = 4037.68
$ws.Range("K136").Value = 9164.7999
$ws.Range("L136").Value = 12113.04
$ws.Range("M136").Value = -6614.7999
$ws.Range("N136").Value = -17213.04

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1145
$ws.Range("I3").Value = 1117.5
$ws.Range("J3").Value = 1200
$ws.Range("K3").Value = 1117.5
$ws.Range("L3").Value = 1200
$ws.Range("M3").Value = -1003.5
$ws.Range("N3").Value = -1428
$ws.Range("H40").Value = 29448
$ws.Range("J40").Value = 29448
$ws.Range("L40").Value = 29448
$ws.Range("N40").Value = -29978
$ws.Range("H86").Value = 6777.75
$ws.Range("J86").Value = 6359
$ws.Range("L86").Value = 6359
$ws.Range("N86").Value = -8605
$ws.Range("H89").Value = 6777.75
$ws.Range("J89").Value = 6359
$ws.Range("L89").Value = 31795
$ws.Range("N89").Value = -43027
$ws.Range("H96").Value = 11332
$ws.Range("H107").Value = 2103.3872
$ws.Range("I107").Value = 1843.037
$ws.Range("J107").Value = 3860.75
$ws.Range("K107").Value = 1843.037
$ws.Range("L107").Value = 3860.75
$ws.Range("M107").Value = 76.96299999999997
$ws.Range("N107").Value = -7700.75
$ws.Range("H134").Value = 28209.256
$ws.Range("I134").Value = 45233.25
$ws.Range("J134").Value = 6705.263
$ws.Range("K134").Value = 135699.75
$ws.Range("L134").Value = 20115.789
$ws.Range("M134").Value = -133164.75
$ws.Range("N134").Value = -25185.789

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1487.3636
$ws.Range("I16").Value = 1662.3334
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 1662.3334
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = -1375.3334
$ws.Range("N16").Value = -1274
$ws.Range("H17").Value = 2387.5
$ws.Range("I17").Value = 600
$ws.Range("J17").Value = 2983.3333
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 2983.3333
$ws.Range("M17").Value = -426
$ws.Range("N17").Value = -3331.3333
$ws.Range("H94").Value = 3981
$ws.Range("I94").Value = 933.6
$ws.Range("J94").Value = 5773.5884
$ws.Range("K94").Value = 933.6
$ws.Range("L94").Value = 5773.5884
$ws.Range("M94").Value = -482.6
$ws.Range("N94").Value = -6675.5884
$ws.Range("H113").Value = 1487.3636
$ws.Range("I113").Value = 1662.3334
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 1662.3334
$ws.Range("L113").Value = 700
$ws.Range("M113").Value = 507.6666
$ws.Range("N113").Value = -5040

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 319.8
$ws.Range("I68").Value = 299.66666
$ws.Range("J68").Value = 350
$ws.Range("K68").Value = 898.9999799999999
$ws.Range("L68").Value = 1050
$ws.Range("M68").Value = -87.99997999999994
$ws.Range("N68").Value = -2672
$ws.Range("H71").Value = 319.8
$ws.Range("I71").Value = 299.66666
$ws.Range("J71").Value = 350
$ws.Range("K71").Value = 2696.99994
$ws.Range("L71").Value = 3150
$ws.Range("M71").Value = 1359.00006
$ws.Range("N71").Value = -11262
$ws.Range("H131").Value = 889.7593000000001
$ws.Range("J131").Value = 912.1799999999999
$ws.Range("L131").Value = 2736.54
$ws.Range("N131").Value = -12816.54
$ws.Range("H132").Value = 4758.8667
$ws.Range("J132").Value = 13100
$ws.Range("L132").Value = 117900
$ws.Range("N132").Value = -122960

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7046.846
$ws.Range("I70").Value = 4189.6665
$ws.Range("J70").Value = 41333
$ws.Range("K70").Value = 4189.6665
$ws.Range("L70").Value = 41333
$ws.Range("M70").Value = -3919.6665
$ws.Range("N70").Value = -41873
$ws.Range("H73").Value = 7046.846
$ws.Range("I73").Value = 4189.6665
$ws.Range("J73").Value = 41333
$ws.Range("K73").Value = 4189.6665
$ws.Range("L73").Value = 41333
$ws.Range("M73").Value = -3253.6665
$ws.Range("N73").Value = -43205
$ws.Range("H80").Value = 4507.5757
$ws.Range("I80").Value = 4810.7144
$ws.Range("K80").Value = 4810.7144
$ws.Range("M80").Value = -3812.7144
$ws.Range("H83").Value = 4507.5757
$ws.Range("I83").Value = 4810.7144
$ws.Range("K83").Value = 24053.572
$ws.Range("M83").Value = -19061.572
$ws.Range("H97").Value = 2017.8889
$ws.Range("J97").Value = 3040.25
$ws.Range("L97").Value = 3040.25
$ws.Range("N97").Value = -4032.25

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1411.027
$ws.Range("I16").Value = 1356.8667
$ws.Range("J16").Value = 1643.1428
$ws.Range("K16").Value = 1356.8667
$ws.Range("L16").Value = 1643.1428
$ws.Range("M16").Value = -1186.8667
$ws.Range("N16").Value = -1983.1428
$ws.Range("H68").Value = 2238.9656
$ws.Range("I68").Value = 1811.3334
$ws.Range("J68").Value = 2697.1428
$ws.Range("K68").Value = 1811.3334
$ws.Range("L68").Value = 2697.1428
$ws.Range("M68").Value = -1062.3334
$ws.Range("N68").Value = -4195.1428
$ws.Range("H71").Value = 2238.9656
$ws.Range("I71").Value = 1811.3334
$ws.Range("J71").Value = 2697.1428
$ws.Range("K71").Value = 9056.666999999999
$ws.Range("L71").Value = 13485.714
$ws.Range("M71").Value = -5312.666999999999
$ws.Range("N71").Value = -20973.714
$ws.Range("H82").Value = 2309.4167
$ws.Range("I82").Value = 1601.8572
$ws.Range("J82").Value = 3300
$ws.Range("K82").Value = 1601.8572
$ws.Range("L82").Value = 3300
$ws.Range("M82").Value = -1240.8572
$ws.Range("N82").Value = -4022
$ws.Range("H85").Value = 2309.4167
$ws.Range("I85").Value = 1601.8572
$ws.Range("J85").Value = 3300
$ws.Range("K85").Value = 1601.8572
$ws.Range("L85").Value = 3300
$ws.Range("M85").Value = -353.8571999999999
$ws.Range("N85").Value = -5796
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H136").Value = 5563.057
$ws.Range("I136").Value = 3179.1052
$ws.Range("K136").Value = 9537.3156
$ws.Range("M136").Value = -6987.3156

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()
$ws.Range("H132").Value = 27616.375
$ws.Range("I132").Value = 85270.75
$ws.Range("J132").Value = 2907.3572
$ws.Range("K132").Value = 255812.25
$ws.Range("L132").Value = 8722.071599999999
$ws.Range("M132").Value = -253282.25
$ws.Range("N132").Value = -13782.0716
